$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Course" column (C) describing which academic year/course each
# group in column A belongs to (or "Магистры" for master's groups).
$ws.Range("C1").Value = "Course"

for ($r = 2; $r -le 14; $r++) {
    $ws.Cells.Item($r, 3).Value = 1
}
for ($r = 15; $r -le 27; $r++) {
    $ws.Cells.Item($r, 3).Value = 2
}
for ($r = 28; $r -le 38; $r++) {
    $ws.Cells.Item($r, 3).Value = 3
}
for ($r = 39; $r -le 48; $r++) {
    $ws.Cells.Item($r, 3).Value = 4
}
for ($r = 49; $r -le 57; $r++) {
    $ws.Cells.Item($r, 3).Value = "Магистры"
}

# Scroll / selection state matching the authored workbook.
[void]$ws.Range("C49:C57").Select()
$excel.ActiveWindow.ScrollRow = 39
